$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 27
$prev = $row - 1

# Values
$ws.Cells.Item($row, 1).Value = 26
$ws.Cells.Item($row, 2).Value = "israel"
$ws.Cells.Item($row, 3).Value = "ligat-ha-al"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45195.8125
$ws.Cells.Item($row, 6).Value = "Hapoel Tel Aviv"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "Maccabi Tel Aviv"
$ws.Cells.Item($row, 9).Value = 5
$ws.Cells.Item($row, 10).Value = 5.57
$ws.Cells.Item($row, 11).Value = "19/09/2023 18:42"
$ws.Cells.Item($row, 12).Value = 7.87
$ws.Cells.Item($row, 13).Value = "26/09/2023 19:27"
$ws.Cells.Item($row, 14).Value = 4.11
$ws.Cells.Item($row, 15).Value = "19/09/2023 18:42"
$ws.Cells.Item($row, 16).Value = 4.61
$ws.Cells.Item($row, 17).Value = "26/09/2023 19:27"
$ws.Cells.Item($row, 18).Value = 1.57
$ws.Cells.Item($row, 19).Value = "19/09/2023 18:42"
$ws.Cells.Item($row, 20).Value = 1.41
$ws.Cells.Item($row, 21).Value = "26/09/2023 19:23"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-tel-aviv-maccabi-tel-aviv/UajxMGco/"

# Copy cell formatting (styles) from the row above for the styled columns
# (column A uses the bold/bordered "index" style, column E uses the
# date-time number format style), matching the existing per-row pattern.
$ws.Range("A" + $prev).Copy()
$ws.Range("A" + $row).PasteSpecial(-4122)

$ws.Range("E" + $prev).Copy()
$ws.Range("E" + $row).PasteSpecial(-4122)

$excel.CutCopyMode = 0
